# Replace the old "Project Clear-LucidHealth-FVA-105379" project entry with
# "Bartush-Cotton Creek Capital Management LLC-FVA-110095" across every
# worksheet that references it (A2 on Project_Title, RateSheetManagement,
# and WeeklyEntryMatrix).
$wb = $excel.ActiveWorkbook

$newValue = "Bartush-Cotton Creek Capital Management LLC-FVA-110095"
$sheetNames = @("Project_Title", "RateSheetManagement", "WeeklyEntryMatrix")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A2").Value = $newValue
}
